$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.374.13"
$ws.Range("E2").Value = "  +0.47%  "

$ws.Range("D3").Value = "1.873.86"
$ws.Range("E3").Value = "  +0.67%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7111"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3113"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07780"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.67%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08461"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").Value = "1.868.53"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.233"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7121"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.31%  "

$ws.Range("D16").Value = "29.376.95"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.046"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008216"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +5.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.76%  "

$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").Value = "2.119.05"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.777"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1600"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.49"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.063"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.94%  "

$ws.Range("E28").Value = "  -0.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.511"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.79%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.431"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.14%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.287"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.24%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.305"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05275"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.14%  "

$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("E35").Value = "  +0.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7461"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -9.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.697"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("E38").Value = "  +0.78%  "

$ws.Range("D39").Value = "1.221.30"
$ws.Range("E39").Value = "  +5.11%  "

$ws.Range("E40").Value = "  +1.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.486"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.68%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8908"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.06%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "2.017.00"
$ws.Range("E46").Value = "  -0.63%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.810"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.68%  "

$ws.Range("E48").Value = "  +0.52%  "

$ws.Range("E49").Value = "  +6.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.375"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.12%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4323"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.05%  "
